# Applies a cyclic re-ordering of the observation data (columns A,B,E,F,G,H,Q,R,AI)
# across rows 47-51 of the active worksheet, matching the upstream re-sequencing
# of the "Id" (column A) values. Row-local metadata (C,D,I,J,P,S,T,U,V,W,Y,AA,
# AD,AE,AG,AW,AX,AY) is identical for every row in this block, so it is left
# untouched; only the species-observation fields that travel with each "Id" move.

$ws = $excel.ActiveWorkbook.ActiveSheet

# Snapshot the "before" values for the columns that move together, keyed by row.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AI")
$rows = 47..51

$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# after_row -> source before_row (the 5-cycle observed in the data).
$mapping = @{
    47 = 49
    48 = 47
    49 = 51
    50 = 48
    51 = 50
}

foreach ($r in $rows) {
    $src = $mapping[$r]
    $srcVals = $snapshot[$src]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $srcVals[$col]
    }
}
